$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new user record as row 33 (id 110032 / Ewan Marsh)
$row = 33
$ws.Cells.Item($row, 1).Value = 110032
$ws.Cells.Item($row, 2).Value = 9317596770
$ws.Cells.Item($row, 3).Value = "Ewan Marsh"
$ws.Cells.Item($row, 4).Value = "ewan.marsh@xyz.com"
$ws.Cells.Item($row, 5).Value = 818876433
$ws.Cells.Item($row, 6).Value = "ACT"
$ws.Cells.Item($row, 7).Value = "eng"
$ws.Cells.Item($row, 8).Value = "PWD"
$ws.Cells.Item($row, 9).Value = $true
$ws.Cells.Item($row, 10).Value = "superadmin"
$ws.Cells.Item($row, 11).Value = "now()"
$ws.Cells.Item($row, 12).Value = "now()"

# Reproduce the final UI selection state: entire columns M:XFD selected,
# scrolled back to the top of the sheet.
$r = $ws.Range($ws.Columns.Item(13), $ws.Columns.Item(16384))
$r.Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
